# Corrected some missing information: the "details" text for the last
# several experience rows had been entered one column too far to the
# right (column H instead of column G). Move that data back into column
# G and clear out the now-empty column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("experience")

for ($r = 13; $r -le 21; $r++) {
    $details = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 7).Value2 = $details
    $ws.Cells.Item($r, 8).ClearContents()
}

# G3 had picked up an extra (unused) style along the way; reset it back
# to the default "Normal" style used by its neighbours in column G.
$ws.Cells.Item(3, 7).Style = "Normal"

# Reflect the sheet the author was actually working in when the fix was
# made: "experience" becomes the active/selected sheet, with the newly
# corrected range selected.
[void]$ws.Activate()
$ws.Range("G13:G21").Select() | Out-Null
